$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E contain text-formatted numbers/percentages (e.g. "62.428.57",
# "  -2.64%  "). Force text number format before assigning so Excel does not
# auto-coerce these strings into numeric values, then restore the default
# "Normal" style (matching the source workbook, which has no explicit style on
# these cells) once all the writes are done.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.428.57"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "3.005.18"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "584.44"
$ws.Range("D6").Value = "146.59"
$ws.Range("E6").Value = "  -6.59%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").Value = "3.004.49"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("E10").Value = "  -5.91%  "
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  -4.84%  "
$ws.Range("D14").Value = "34.71"
$ws.Range("E14").Value = "  -6.28%  "
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "3.499.02"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "62.413.73"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "3.004.62"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "459.12"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").Value = "13.90"
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("D24").Value = "80.08"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  -9.01%  "
$ws.Range("D26").Value = "12.22"
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  -6.48%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  -5.27%  "
$ws.Range("D33").Value = "26.93"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("D35").Value = "1.03"
$ws.Range("E35").Value = "  -3.61%  "
$ws.Range("D36").Value = "0.0₃0789"
$ws.Range("E36").Value = "  -6.64%  "
$ws.Range("D37").Value = "5.74"
$ws.Range("E37").Value = "  -5.04%  "
$ws.Range("E38").Value = "  -6.80%  "
$ws.Range("D39").Value = "50.02"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "8.95"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  -10.48%  "
$ws.Range("D42").Value = "409.60"
$ws.Range("E42").Value = "  -7.69%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  -5.29%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0353"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("D46").Value = "2.771.94"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "39.09"
$ws.Range("D48").Value = "127.49"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").Value = "23.73"
$ws.Range("E51").Value = "  -9.69%  "

$ws.Range("D2:E51").Style = "Normal"
